$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B8").Value = 5.634499999999997
$ws.Range("B10").Value = 5.900099999999997
$ws.Range("B12").Value = 4.683999999999999
$ws.Range("B18").Value = 6.765199999999992
$ws.Range("B25").Value = 6.000299999999997
$ws.Range("B37").Value = 8.749800000000004
$ws.Range("B55").Value = 6.159599999999992
$ws.Range("B68").Value = 4.688299999999995
$ws.Range("B77").Value = 8.784700000000001
$ws.Range("B78").Value = 9.267600000000002
$ws.Range("B79").Value = 8.675700000000004
$ws.Range("B80").Value = 9.282200000000001
$ws.Range("B81").Value = 5.537700000000005
$ws.Range("B82").Value = 5.284900000000001
$ws.Range("B84").Value = 5.639500000000001
$ws.Range("B101").Value = 5.684699999999995
$ws.Range("B102").Value = 7.870599999999997
